$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.135.29"
$ws.Range("E2").Value = "  -3.89%  "
$ws.Range("D3").Value = "2.235.78"
$ws.Range("E3").Value = "  -4.78%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.69"
$ws.Range("E5").Value = "  -3.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.630"
$ws.Range("E6").Value = "  -6.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.68"
$ws.Range("E7").Value = "  -3.57%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.557"
$ws.Range("E9").Value = "  -7.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0987"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.07"
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "35.12"
$ws.Range("E12").Value = "  +7.25%  "
$ws.Range("E13").Value = "  -2.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.74"
$ws.Range("E14").Value = "  -7.94%  "
$ws.Range("D15").Value = "2.571.01"
$ws.Range("E15").Value = "  -4.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.72"
$ws.Range("E16").Value = "  -10.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.865"
$ws.Range("E17").Value = "  -4.22%  "
$ws.Range("D18").Value = "2.246.54"
$ws.Range("E18").Value = "  -4.41%  "
$ws.Range("D19").Value = "42.038.97"
$ws.Range("E19").Value = "  -3.77%  "
$ws.Range("D20").Value = "0.0₃0982"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.00"
$ws.Range("E21").Value = "  -5.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.16"
$ws.Range("E22").Value = "  -8.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.97"
$ws.Range("E23").Value = "  -8.82%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.63"
$ws.Range("E25").Value = "  -2.91%  "
$ws.Range("E26").Value = "  -8.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.36"
$ws.Range("E27").Value = "  -5.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.96"
$ws.Range("E28").Value = "  -6.32%  "
$ws.Range("E29").Value = "  -8.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.45"
$ws.Range("E30").Value = "  -6.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.52"
$ws.Range("E31").Value = "  -9.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.126"
$ws.Range("E33").Value = "  -8.43%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.26"
$ws.Range("E34").Value = "  -4.34%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0705"
$ws.Range("E35").Value = "  -7.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.76"
$ws.Range("E36").Value = "  -8.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.54"
$ws.Range("E37").Value = "  -7.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.04"
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.22"
$ws.Range("E39").Value = "  -6.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "20.68"
$ws.Range("E40").Value = "  +8.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0264"
$ws.Range("E41").Value = "  -6.17%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.07"
$ws.Range("E42").Value = "  -7.57%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.11"
$ws.Range("E43").Value = "  +6.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.73"
$ws.Range("E44").Value = "  -4.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0995"
$ws.Range("E45").Value = "  -10.66%  "
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.185"
$ws.Range("E47").Value = "  -8.72%  "
$ws.Range("B48").Value = "SynthetixNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.34"
$ws.Range("E48").Value = "  +7.18%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.16"
$ws.Range("E49").Value = "  -6.10%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.32"
$ws.Range("E50").Value = "  -7.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.80"
$ws.Range("E51").Value = "  -0.95%  "
